# "gave up on individual position" - collapse a few individually-tracked
# schedule positions ("1" / "2") into the shared "2*" marker, matching the
# other cells in the same rows that already use "2*".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C9").Value = "2*"

$ws.Range("E11").Value = "2*"
$ws.Range("G11").Value = "2*"
$ws.Range("I11").Value = "2*"
$ws.Range("K11").Value = "2*"
